$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 160.9723563333334
$ws.Range("H2").Value = 482.917069
$ws.Range("I2").Value = 0.3931645655589854
$ws.Range("J2").Value = 0.3931645655589854
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 13.69681033333333
$ws.Range("N2").Value = 41.090431
$ws.Range("O2").Value = 0.1107101339353595
$ws.Range("P2").Value = 0.1107101339353595
$ws.Range("Q2").Value = 2204.807833607415
$ws.Range("R2").Value = 19843.27050246674
$ws.Range("S2").Value = 0.04352730171167271
$ws.Range("T2").Value = 0.04352730171167271
$ws.Range("G3").Value = 160.9723563333334
$ws.Range("H3").Value = 482.917069
$ws.Range("I3").Value = 0.3931645655589854
$ws.Range("J3").Value = 0.3931645655589854
$ws.Range("O3").Value = 0.8037307792188669
$ws.Range("P3").Value = 0.803730779218867
$ws.Range("Q3").Value = 16006.41111289606
$ws.Range("R3").Value = 144057.7000160646
$ws.Range("S3").Value = 0.3159984626379706
$ws.Range("T3").Value = 0.3159984626379707
$ws.Range("G4").Value = 160.9723563333334
$ws.Range("H4").Value = 482.917069
$ws.Range("I4").Value = 0.3931645655589854
$ws.Range("J4").Value = 0.3931645655589854
$ws.Range("O4").Value = 0.08555908684577355
$ws.Range("P4").Value = 0.08555908684577354
$ws.Range("Q4").Value = 1703.921205840117
$ws.Range("R4").Value = 15335.29085256106
$ws.Range("S4").Value = 0.03363880120934206
$ws.Range("T4").Value = 0.03363880120934205
$ws.Range("I5").Value = 0.2197635343237224
$ws.Range("J5").Value = 0.2197635343237224
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 13.69681033333333
$ws.Range("N5").Value = 41.090431
$ws.Range("O5").Value = 0.1107101339353595
$ws.Range("P5").Value = 0.1107101339353595
$ws.Range("Q5").Value = 1232.400893832589
$ws.Range("R5").Value = 11091.6080444933
$ws.Range("S5").Value = 0.02433005031908728
$ws.Range("T5").Value = 0.02433005031908728
$ws.Range("I6").Value = 0.2197635343237224
$ws.Range("J6").Value = 0.2197635343237224
$ws.Range("O6").Value = 0.8037307792188669
$ws.Range("P6").Value = 0.803730779218867
$ws.Range("S6").Value = 0.1766307166858976
$ws.Range("T6").Value = 0.1766307166858976
$ws.Range("I7").Value = 0.2197635343237224
$ws.Range("J7").Value = 0.2197635343237224
$ws.Range("O7").Value = 0.08555908684577355
$ws.Range("P7").Value = 0.08555908684577354
$ws.Range("Q7").Value = 952.4249619804151
$ws.Range("R7").Value = 8571.824657823736
$ws.Range("S7").Value = 0.0188027673187375
$ws.Range("T7").Value = 0.0188027673187375
$ws.Range("I8").Value = 0.3870719001172923
$ws.Range("J8").Value = 0.3870719001172923
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 13.69681033333333
$ws.Range("N8").Value = 41.090431
$ws.Range("O8").Value = 0.1107101339353595
$ws.Range("P8").Value = 0.1107101339353595
$ws.Range("Q8").Value = 2170.641080878071
$ws.Range("R8").Value = 19535.76972790263
$ws.Range("S8").Value = 0.04285278190459953
$ws.Range("T8").Value = 0.04285278190459953
$ws.Range("I9").Value = 0.3870719001172923
$ws.Range("J9").Value = 0.3870719001172923
$ws.Range("O9").Value = 0.8037307792188669
$ws.Range("P9").Value = 0.803730779218867
$ws.Range("S9").Value = 0.3111015998949987
$ws.Range("T9").Value = 0.3111015998949988
$ws.Range("I10").Value = 0.3870719001172923
$ws.Range("J10").Value = 0.3870719001172923
$ws.Range("O10").Value = 0.08555908684577355
$ws.Range("P10").Value = 0.08555908684577354
$ws.Range("R10").Value = 15097.64787859958
$ws.Range("S10").Value = 0.03311751831769399
$ws.Range("T10").Value = 0.03311751831769399
